$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay as text (they use dot as thousands separator)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "33.671.89"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "1.772.83"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "224.23"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "0.545"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "31.94"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").Value = "0.0683"
$ws.Range("E10").Value = "  -4.48%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").Value = "2.027.56"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "11.08"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").Value = "1.777.46"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "33.695.66"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "66.49"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").Value = "0.0₃0773"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").Value = "238.06"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "10.57"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").Value = "159.61"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "16.10"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").Value = "7.00"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").Value = "0.0510"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").Value = "1.380.79"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "77.95"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("E42").Value = "  -4.06%  "
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").Value = "13.53"
$ws.Range("E44").Value = "  +13.92%  "
$ws.Range("D45").Value = "1.08"
$ws.Range("E45").Value = "  +4.08%  "
$ws.Range("D46").Value = "0.0500"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("E47").Value = "  +12.57%  "
$ws.Range("D48").Value = "107.27"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").Value = "5.83"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "1.927.91"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "0.999"
